$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 91 and 92: the two match records were swapped (B,F,G,H,I,J,K..AC)
# Row 91 becomes what used to be row 92's data, and vice versa.
# ---------------------------------------------------------------------------

# New row 91 (previously row 92's content)
$ws.Range("B91").Value  = 6924568
$ws.Range("F91").Value  = "Atletico Morelia"
$ws.Range("G91").Value  = "Atlante"
$ws.Range("H91").Value  = 0
$ws.Range("I91").Value  = 1
$ws.Range("J91").Value  = "A"
$ws.Range("K91").Value  = 2.4
$ws.Range("L91").Value  = 3
$ws.Range("M91").Value  = 2.875
$ws.Range("N91").Value  = 2.7
$ws.Range("O91").Value  = 3.1
$ws.Range("P91").Value  = 2.8
$ws.Range("Q91").Value  = 0
$ws.Range("R91").Value  = 1.85
$ws.Range("S91").Value  = 1.95
$ws.Range("T91").Value  = 2.25
$ws.Range("U91").Value  = 1.975
$ws.Range("V91").Value  = 1.725
$ws.Range("W91").Value  = -1
$ws.Range("X91").Value  = -1
$ws.Range("Y91").Value  = 1.8
$ws.Range("Z91").Value  = -1
$ws.Range("AA91").Value = 0.95
$ws.Range("AB91").Value = -1
$ws.Range("AC91").Value = 0.7250000000000001

# New row 92 (previously row 91's content)
$ws.Range("B92").Value  = 6924569
$ws.Range("F92").Value  = "Venados FC"
$ws.Range("G92").Value  = "Dorados"
$ws.Range("H92").Value  = 4
$ws.Range("I92").Value  = 1
$ws.Range("J92").Value  = "H"
$ws.Range("K92").Value  = 1.615
$ws.Range("L92").Value  = 4
$ws.Range("M92").Value  = 4.5
$ws.Range("N92").Value  = 1.5
$ws.Range("O92").Value  = 4.75
$ws.Range("P92").Value  = 5.75
$ws.Range("Q92").Value  = -1.25
$ws.Range("R92").Value  = 1.925
$ws.Range("S92").Value  = 1.875
$ws.Range("T92").Value  = 3
$ws.Range("U92").Value  = 1.75
$ws.Range("V92").Value  = 1.95
$ws.Range("W92").Value  = 0.5
$ws.Range("X92").Value  = -1
$ws.Range("Y92").Value  = -1
$ws.Range("Z92").Value  = 0.925
$ws.Range("AA92").Value = -1
$ws.Range("AB92").Value = 0.75
$ws.Range("AC92").Value = -1

# ---------------------------------------------------------------------------
# Row 173: match result came in, add FTHG/FTAG/FTR and PL_* columns, refresh
# closing odds.
# ---------------------------------------------------------------------------
$ws.Range("H173").Value = 1
$ws.Range("I173").Value = 0
$ws.Range("J173").Value = "H"
$ws.Range("N173").Value = 2.8
$ws.Range("O173").Value = 3.2
$ws.Range("R173").Value = 1.95
$ws.Range("S173").Value = 1.85
$ws.Range("T173").Value = 2
$ws.Range("U173").Value = 1.825
$ws.Range("V173").Value = 1.975
$ws.Range("W173").Value = 1.8
$ws.Range("X173").Value = -1
$ws.Range("Y173").Value = -1
$ws.Range("Z173").Value = 0.95
$ws.Range("AA173").Value = -1
$ws.Range("AB173").Value = -1
$ws.Range("AC173").Value = 0.9750000000000001

# ---------------------------------------------------------------------------
# Row 174: match result came in, add FTHG/FTAG/FTR and PL_* columns, refresh
# closing odds.
# ---------------------------------------------------------------------------
$ws.Range("H174").Value = 1
$ws.Range("I174").Value = 1
$ws.Range("J174").Value = "D"
$ws.Range("N174").Value = 1.727
$ws.Range("P174").Value = 5
$ws.Range("R174").Value = 1.95
$ws.Range("S174").Value = 1.85
$ws.Range("U174").Value = 1.85
$ws.Range("V174").Value = 1.95
$ws.Range("W174").Value = -1
$ws.Range("X174").Value = 2.75
$ws.Range("Y174").Value = -1
$ws.Range("Z174").Value = -1
$ws.Range("AA174").Value = 0.8500000000000001
$ws.Range("AB174").Value = -1
$ws.Range("AC174").Value = 0.95

# ---------------------------------------------------------------------------
# Row 175: match result came in, add FTHG/FTAG/FTR and PL_* columns, refresh
# closing odds.
# ---------------------------------------------------------------------------
$ws.Range("H175").Value = 1
$ws.Range("I175").Value = 1
$ws.Range("J175").Value = "D"
$ws.Range("N175").Value = 3.3
$ws.Range("O175").Value = 3.75
$ws.Range("P175").Value = 2.05
$ws.Range("Q175").Value = 0.5
$ws.Range("R175").Value = 1.8
$ws.Range("S175").Value = 2
$ws.Range("T175").Value = 2.75
$ws.Range("U175").Value = 1.925
$ws.Range("V175").Value = 1.875
$ws.Range("W175").Value = -1
$ws.Range("X175").Value = 2.75
$ws.Range("Y175").Value = -1
$ws.Range("Z175").Value = 0.8
$ws.Range("AA175").Value = -1
$ws.Range("AB175").Value = -1
$ws.Range("AC175").Value = 0.875

# ---------------------------------------------------------------------------
# Row 176: becomes the fixture that used to be row 177 (new odds pulled in),
# W..AA remain 0 (match not yet played).
# ---------------------------------------------------------------------------
$ws.Range("B176").Value = 7641684
$ws.Range("E176").Value = 45351.92013888889
$ws.Range("F176").Value = "Tepatitlan FC"
$ws.Range("G176").Value = "Club Celaya"
$ws.Range("K176").Value = 4
$ws.Range("L176").Value = 3.3
$ws.Range("M176").Value = 1.95
$ws.Range("N176").Value = 4.75
$ws.Range("O176").Value = 3.5
$ws.Range("P176").Value = 1.833
$ws.Range("Q176").Value = 0.75
$ws.Range("R176").Value = 1.8
$ws.Range("S176").Value = 2
$ws.Range("T176").Value = 2.25
$ws.Range("U176").Value = 1.95
$ws.Range("V176").Value = 1.85

# ---------------------------------------------------------------------------
# Row 177: becomes the fixture that used to be row 178 (new odds pulled in),
# W..AA remain 0 (match not yet played).
# ---------------------------------------------------------------------------
$ws.Range("B177").Value = 7641685
$ws.Range("E177").Value = 45352.00347222222
$ws.Range("F177").Value = "Atletico Morelia"
$ws.Range("G177").Value = "Cimarrones de Sonora FC"
$ws.Range("K177").Value = 1.8
$ws.Range("L177").Value = 3.25
$ws.Range("M177").Value = 4
$ws.Range("N177").Value = 1.65
$ws.Range("O177").Value = 4
$ws.Range("P177").Value = 5.25
$ws.Range("Q177").Value = -0.75
$ws.Range("R177").Value = 1.85
$ws.Range("S177").Value = 1.95
$ws.Range("T177").Value = 2.5
$ws.Range("U177").Value = 1.875
$ws.Range("V177").Value = 1.925

# ---------------------------------------------------------------------------
# Row 178: a brand-new fixture is pulled in (new id/date/teams/odds),
# W..AA remain 0 (match not yet played).
# ---------------------------------------------------------------------------
$ws.Range("B178").Value = 7701489
$ws.Range("E178").Value = 45352.83333333334
$ws.Range("F178").Value = "Correcaminos"
$ws.Range("G178").Value = "Tlaxcala FC"
$ws.Range("K178").Value = 1.571
$ws.Range("L178").Value = 3.6
$ws.Range("M178").Value = 5
$ws.Range("N178").Value = 1.7
$ws.Range("O178").Value = 3.6
$ws.Range("P178").Value = 5.75
$ws.Range("Q178").Value = -0.75
$ws.Range("R178").Value = 1.825
$ws.Range("S178").Value = 1.975
$ws.Range("T178").Value = 2.25
$ws.Range("U178").Value = 1.95
$ws.Range("V178").Value = 1.85

# ---------------------------------------------------------------------------
# Row 179: new row appended at the bottom of the sheet for a further new
# fixture. Copy formatting from row 178 (style "1" on column A, date style
# "2" on column E) before filling in the values.
# ---------------------------------------------------------------------------
$ws.Range("A178").Copy()
$ws.Range("A179").PasteSpecial(-4122)
$ws.Range("E178").Copy()
$ws.Range("E179").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A179").Value  = 177
$ws.Range("B179").Value  = 7641686
$ws.Range("C179").Value  = "Mexico Liga de Expansion"
$ws.Range("D179").Value  = "Mexico Liga de Expansion"
$ws.Range("E179").Value  = 45353.83680555555
$ws.Range("F179").Value  = "Cancun FC"
$ws.Range("G179").Value  = "Club Atletico La Paz"
$ws.Range("K179").Value  = 2
$ws.Range("L179").Value  = 3.1
$ws.Range("M179").Value  = 3.4
$ws.Range("N179").Value  = 1.8
$ws.Range("O179").Value  = 3.4
$ws.Range("P179").Value  = 5
$ws.Range("Q179").Value  = -0.5
$ws.Range("R179").Value  = 1.8
$ws.Range("S179").Value  = 2
$ws.Range("T179").Value  = 2.25
$ws.Range("U179").Value  = 1.775
$ws.Range("V179").Value  = 2.025
$ws.Range("W179").Value  = 0
$ws.Range("X179").Value  = 0
$ws.Range("Y179").Value  = 0
$ws.Range("Z179").Value  = 0
$ws.Range("AA179").Value = 0
